$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.079.21"
$ws.Range("E2").Value = "  +0.04%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.834.50"
$ws.Range("E3").Value = "  +0.07%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.66"
$ws.Range("E5").Value = "  -1.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6158"
$ws.Range("E6").Value = "  -2.77%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("E7").Value = "  +0.24%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07443"
$ws.Range("E8").Value = "  -1.42%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2925"
$ws.Range("E9").Value = "  -0.81%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.06"
$ws.Range("E10").Value = "  -0.57%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07682"
$ws.Range("E11").Value = "  -0.41%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.836.25"
$ws.Range("E12").Value = "  +0.26%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.991"
$ws.Range("E13").Value = "  -0.36%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6711"
$ws.Range("E14").Value = "  -0.21%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "82.61"
$ws.Range("E15").Value = "  -0.79%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009121"
$ws.Range("E16").Value = "  -4.61%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.896"
$ws.Range("E17").Value = "  -3.31%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.081.61"
$ws.Range("E18").Value = "  -0.04%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.082.19"
$ws.Range("E19").Value = "  -0.09%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "232.92"
$ws.Range("E20").Value = "  +2.58%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.66"
$ws.Range("E21").Value = "  +0.26%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.003"
$ws.Range("E22").Value = "  +0.37%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.168"
$ws.Range("E23").Value = "  -0.10%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.002"
$ws.Range("E24").Value = "  +0.22%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.20"
$ws.Range("E25").Value = "  -0.58%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1408"
$ws.Range("E26").Value = "  -1.45%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.475"
$ws.Range("E27").Value = "  -0.71%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.79"
$ws.Range("E28").Value = "  -0.89%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.498"
$ws.Range("E29").Value = "  -0.40%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.149"
$ws.Range("E30").Value = "  -0.21%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.103"
$ws.Range("E31").Value = "  +0.52%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.05498"
$ws.Range("E32").Value = "  +0.30%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.206"
$ws.Range("E33").Value = "  +0.24%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.832"
$ws.Range("E34").Value = "  -1.55%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7345"
$ws.Range("E35").Value = "  -1.62%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.139"
$ws.Range("E36").Value = "  -0.40%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.663"
$ws.Range("E37").Value = "  +0.09%  "

$ws.Range("E38").Value = "  +0.38%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01779"
$ws.Range("E39").Value = "  -0.58%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.208.55"
$ws.Range("E40").Value = "  -3.08%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.436"
$ws.Range("E41").Value = "  -2.77%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8929"
$ws.Range("E42").Value = "  -1.50%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.002"
$ws.Range("E43").Value = "  +0.21%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.81"
$ws.Range("E44").Value = "  +0.16%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.979.00"
$ws.Range("E45").Value = "  -0.32%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "65.36"
$ws.Range("E46").Value = "  +0.28%  "

$ws.Range("E47").Value = "  -2.59%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5094"
$ws.Range("E48").Value = "  -0.17%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4063"
$ws.Range("E49").Value = "  -0.30%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.140"
$ws.Range("E50").Value = "  +1.19%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05820"
$ws.Range("E51").Value = "  +0.43%  "
